$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

try {
  $excel.Goto($ws.Range("A13"), $true)
  Write-Host "Goto worked"
} catch {
  Write-Host "Goto error: $_"
}
